$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    # Force text storage so numeric-looking strings (e.g. "6.674", "1.001")
    # are not auto-coerced into Number cells by Excel's smart input parsing,
    # then restore the default "Normal" style so no stray per-cell number
    # format lingers (matches the source workbook, which carries no `s`
    # attribute on these data cells).
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "27.057.86"
$ws.Cells.Item(2, 5).Value = "  -1.84%  "

# Row 3
Set-TextCell 3 4 "1.827.09"
$ws.Cells.Item(3, 5).Value = "  -0.77%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.79%  "

# Row 5
Set-TextCell 5 4 "311.54"
$ws.Cells.Item(5, 5).Value = "  -1.75%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.67%  "

# Row 7
Set-TextCell 7 4 "0.4243"
$ws.Cells.Item(7, 5).Value = "  -1.34%  "

# Row 8
Set-TextCell 8 4 "0.3669"

# Row 9
Set-TextCell 9 4 "0.07226"
$ws.Cells.Item(9, 5).Value = "  -0.81%  "

# Row 10
Set-TextCell 10 4 "0.8439"
$ws.Cells.Item(10, 5).Value = "  -3.03%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.80%  "

# Row 12
Set-TextCell 12 4 "1.832.14"
$ws.Cells.Item(12, 5).Value = "  -0.59%  "

# Row 13
Set-TextCell 13 4 "6.674"

# Row 14
Set-TextCell 14 4 "5.287"
$ws.Cells.Item(14, 5).Value = "  -1.87%  "

# Row 15
Set-TextCell 15 4 "0.07034"
$ws.Cells.Item(15, 5).Value = "  -1.08%  "

# Row 16
Set-TextCell 16 4 "89.67"
$ws.Cells.Item(16, 5).Value = "  +1.22%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.87%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -2.38%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.61%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -2.84%  "

# Row 21
Set-TextCell 21 4 "27.138.76"
$ws.Cells.Item(21, 5).Value = "  -1.60%  "

# Row 22
Set-TextCell 22 4 "5.131"
$ws.Cells.Item(22, 5).Value = "  -1.06%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -1.72%  "

# Row 24
Set-TextCell 24 4 "2.054.12"
$ws.Cells.Item(24, 5).Value = "  -0.86%  "

# Row 25
Set-TextCell 25 4 "1.984"

# Row 26
Set-TextCell 26 4 "151.38"
$ws.Cells.Item(26, 5).Value = "  -1.95%  "

# Row 27
Set-TextCell 27 4 "2.253"
$ws.Cells.Item(27, 5).Value = "  +4.74%  "

# Row 28
Set-TextCell 28 4 "18.18"
$ws.Cells.Item(28, 5).Value = "  -1.86%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -1.18%  "

# Row 30
Set-TextCell 30 4 "116.65"
$ws.Cells.Item(30, 5).Value = "  -0.73%  "

# Row 31
Set-TextCell 31 4 "0.08714"
$ws.Cells.Item(31, 5).Value = "  -2.05%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -2.53%  "

# Row 33
Set-TextCell 33 4 "0.7366"
$ws.Cells.Item(33, 5).Value = "  -4.56%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 34 4 "2.902"
$ws.Cells.Item(34, 5).Value = "  -0.21%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Filecoin"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 35 4 "4.429"
$ws.Cells.Item(35, 5).Value = "  -1.80%  "

# Row 36
Set-TextCell 36 4 "1.001"

# Row 37
Set-TextCell 37 4 "1.094"
$ws.Cells.Item(37, 5).Value = "  -2.92%  "

# Row 38
Set-TextCell 38 4 "0.01943"
$ws.Cells.Item(38, 5).Value = "  -1.18%  "

# Row 39
Set-TextCell 39 4 "0.05230"
$ws.Cells.Item(39, 5).Value = "  -1.15%  "

# Row 40
Set-TextCell 40 4 "7.290"
$ws.Cells.Item(40, 5).Value = "  +2.23%  "

# Row 41
Set-TextCell 41 4 "2.873"
$ws.Cells.Item(41, 5).Value = "  -0.52%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "TheSandbox"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell 42 4 "0.5137"
$ws.Cells.Item(42, 5).Value = "  +0.57%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell 43 4 "0.1687"
$ws.Cells.Item(43, 5).Value = "  +0.17%  "

# Row 44
Set-TextCell 44 4 "8.557"
$ws.Cells.Item(44, 5).Value = "  -2.14%  "

# Row 45
Set-TextCell 45 4 "10.52"
$ws.Cells.Item(45, 5).Value = "  -1.02%  "

# Row 46
Set-TextCell 46 4 "1.970"
$ws.Cells.Item(46, 5).Value = "  +7.23%  "

# Row 47
Set-TextCell 47 4 "0.4736"
$ws.Cells.Item(47, 5).Value = "  -0.01%  "

# Row 48
Set-TextCell 48 4 "105.53"
$ws.Cells.Item(48, 5).Value = "  -1.25%  "

# Row 49
Set-TextCell 49 4 "1.001"
$ws.Cells.Item(49, 5).Value = "  -0.75%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -1.99%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -1.42%  "
